# prelievi.xlsx - "new new 7 maggio" commit
# Adds two new withdrawal rows dated 07/05/2018 (Excel serial 43227) at the
# bottom of the existing table (Sheet1), right after the current last row (12).
#
# Row 13: Bertolotti Daniela | Tessuto lana a quadri | Mt. | 1
# Row 14: Cristina Sarah     | Tela Leggera          | Mt. | 2
#
# Inserting the rows (instead of just writing into the first empty row)
# makes Excel carry the formatting/styles down from the row above, which is
# what the target workbook expects for the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 13 ---
$ws.Rows("13").Insert()
$ws.Range("A13").Value = 43227
$ws.Range("B13").Value = "Bertolotti Daniela"
$ws.Range("C13").Value = "Tessuto lana a quadri"
$ws.Range("D13").Value = "Mt."
$ws.Range("E13").Value = 1

# --- New row 14 ---
$ws.Rows("14").Insert()
$ws.Range("A14").Value = 43227
$ws.Range("B14").Value = "Cristina Sarah"
$ws.Range("C14").Value = "Tela Leggera"
$ws.Range("D14").Value = "Mt."
$ws.Range("E14").Value = 2
